# Adds the "LTSD Parameters" activity-detection threshold table
# (columns K:N, rows 2:5) to the APA_Liste_3 sheet, as produced by the
# "Added excel files for the APA analysis" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / labels (plain text, left-to-right, top-to-bottom so the
#     shared-string table grows in the same order the source workbook used) ---
$ws.Range("K2").Value = "LTSD Parameters"

$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

$ws.Range("K4").Value = "Threshols"
$ws.Range("L4").Value = "Win"
$ws.Range("M4").Value = "Threshold"
$ws.Range("N4").Value = "Win"

# --- Numeric-looking values that must be stored as TEXT (shared strings),
#     not as numbers. Plain `.Value = "200.0"` would get auto-coerced to a
#     number by Excel, so each value is entered as a formula producing the
#     exact string, then flattened to a literal value via copy/paste-special
#     (values only) - this keeps cell type "s" without adding any new
#     number-format/quote-prefix style.
#     N5 is written first so the shared-string insertion order matches the
#     source file (200.0, 7.4, 460.0, 6.0).

$ws.Range("N5").Formula = '="200.0"'
$ws.Range("N5").Copy()
$ws.Range("N5").PasteSpecial(-4163)

$ws.Range("K5").Formula = '="7.4"'
$ws.Range("K5").Copy()
$ws.Range("K5").PasteSpecial(-4163)

$ws.Range("L5").Formula = '="460.0"'
$ws.Range("L5").Copy()
$ws.Range("L5").PasteSpecial(-4163)

$ws.Range("M5").Formula = '="6.0"'
$ws.Range("M5").Copy()
$ws.Range("M5").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# Matches the saved selection state in the edited workbook.
$ws.Range("M5").Select()
